# Actualización automática 2025-12-07 17:44:30
# Update the PRESUPUESTO (column G) values on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G2").Value  = 0
$ws.Range("G3").Value  = 1000
$ws.Range("G4").Value  = 1000
$ws.Range("G6").Value  = 2500
$ws.Range("G8").Value  = 0
$ws.Range("G9").Value  = 1500
$ws.Range("G11").Value = 0
$ws.Range("G12").Value = 12870
$ws.Range("G13").Value = 1000
$ws.Range("G14").Value = 1000
$ws.Range("G16").Value = 500
$ws.Range("G17").Value = 0
$ws.Range("G21").Value = 500
$ws.Range("G22").Value = 1000
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 15000
$ws.Range("G26").Value = 38870
